# Natmi following Dr Hou advice
# Rspo3-Lgr5 L-R pair sheet: expand the single sample row into the full
# sending/target-cluster combination matrix (FAPs/sCs sending clusters x
# FAPs/ECs/M2/sCs target clusters) with recomputed specificity metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T the 16 numeric expression / specificity metrics.
$rowsData = @(
     ,@("FAPs", "Rspo3", "Lgr5", "ECs", 3, 1, 4.83236, 14.49708, 0.975350813525687, 0.975350813525687, 2, 0.6666666666666666, 0.1978116666666667, 0.593435, 0.0447866041160341, 0.04478660411603409, 0.9558971855333336, 8.603074669800002, 0.04368265075962674, 0.04368265075962673)
     ,@("FAPs", "Rspo3", "Lgr5", "FAPs", 3, 1, 4.83236, 14.49708, 0.975350813525687, 0.975350813525687, 3, 1, 1.940645333333333, 5.821936, 0.4393821443306968, 0.4393821443306967, 9.377896882986667, 84.40107194688001, 0.4285517319216059, 0.4285517319216058)
     ,@("FAPs", "Rspo3", "Lgr5", "M2", 3, 1, 4.83236, 14.49708, 0.975350813525687, 0.975350813525687, 1, 0.3333333333333333, 0.046181, 0.138543, 0.0104558553068958, 0.0104558553068958, 0.22316321716, 2.00846895444, 0.01019812697968769, 0.01019812697968769)
     ,@("FAPs", "Rspo3", "Lgr5", "sCs", 3, 1, 4.83236, 14.49708, 0.975350813525687, 0.975350813525687, 2, 0.6666666666666666, 2.232121666666667, 6.696365, 0.5053753962463735, 0.5053753962463734, 10.78641545713333, 97.0777391142, 0.4929183038647668, 0.4929183038647667)
     ,@("sCs", "Rspo3", "Lgr5", "ECs", 2, 0.6666666666666666, 0.122124, 0.366372, 0.02464918647431296, 0.02464918647431296, 2, 0.6666666666666666, 0.1978116666666667, 0.593435, 0.0447866041160341, 0.04478660411603409, 0.02415755198, 0.21741796782, 0.001103953356407357, 0.001103953356407357)
     ,@("sCs", "Rspo3", "Lgr5", "FAPs", 2, 0.6666666666666666, 0.122124, 0.366372, 0.02464918647431296, 0.02464918647431296, 3, 1, 1.940645333333333, 5.821936, 0.4393821443306968, 0.4393821443306967, 0.236999370688, 2.132994336192, 0.01083041240909084, 0.01083041240909084)
     ,@("sCs", "Rspo3", "Lgr5", "M2", 2, 0.6666666666666666, 0.122124, 0.366372, 0.02464918647431296, 0.02464918647431296, 1, 0.3333333333333333, 0.046181, 0.138543, 0.0104558553068958, 0.0104558553068958, 0.005639808444, 0.050758275996, 0.0002577283272081094, 0.0002577283272081094)
     ,@("sCs", "Rspo3", "Lgr5", "sCs", 2, 0.6666666666666666, 0.122124, 0.366372, 0.02464918647431296, 0.02464918647431296, 2, 0.6666666666666666, 2.232121666666667, 6.696365, 0.5053753962463735, 0.5053753962463734, 0.27259562642, 2.45336063778, 0.01245709238160666, 0.01245709238160666)
)

$startRow = 2
for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $rowsData[$i]
    for ($c = 1; $c -le $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
